# Adding info about Jiuchen seminar attendance
# Fill in the Sep 2024 seminar attendance counts for the first two
# "2024 - Fall" sessions (Jiuchen's seminars). Columns G/H hold
# undergrad/grad (or similar) headcounts; I, K and the summary rows
# (17-19) are formulas and will recalculate automatically, as will the
# cross-sheet references on "Attendance Descriptives" and the chart
# caches that are sourced from them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024 - Fall")

$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 5
$ws.Range("G6").Value = 6
$ws.Range("H6").Value = 7

# Match the author's final selection on this sheet after entering the data.
$ws.Activate()
$ws.Range("F8").Select()
